# Rename/reassign income categories, resize the category column, and
# update the active view so the "share of categories" becomes the visual
# focus of the sheet (per commit message: categories are now visualized).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$categoryByRow = @{
    2 = "Private tutoring"
    3 = "Parents transfer"
    4 = "Private tutoring"
    5 = "Parents transfer"
    6 = "Parents transfer"
    7 = "Parents transfer"
    8 = "Scholarship"
    9 = "Private tutoring"
    10 = "Parents transfer"
    11 = "Parents transfer"
    12 = "Parents transfer"
    13 = "Scholarship"
    14 = "Private tutoring"
    15 = "Parents transfer"
    16 = "Parents transfer"
    17 = "Parents transfer"
    18 = "Tax refund"
    19 = "Parents transfer"
    20 = "Parents transfer"
    21 = "Scholarship"
    22 = "Private tutoring"
    23 = "Parents transfer"
    24 = "Private tutoring"
    25 = "Parents transfer"
    26 = "Parents transfer"
    27 = "Parents transfer"
    28 = "Parents transfer"
    29 = "Scholarship"
    30 = "Parents transfer"
    31 = "Parents transfer"
    32 = "Parents transfer"
    33 = "Parents transfer"
    34 = "Parents transfer"
    35 = "Parents transfer"
    36 = "Parents transfer"
    37 = "Parents transfer"
    38 = "Parents transfer"
    39 = "Parents transfer"
    40 = "Scholarship"
    41 = "Parents transfer"
    42 = "Parents transfer"
    43 = "Parents transfer"
    44 = "Parents transfer"
    45 = "Parents transfer"
    46 = "Parents transfer"
    47 = "Parents transfer"
    48 = "Parents transfer"
    49 = "Parents transfer"
    50 = "Parents transfer"
    51 = "Parents transfer"
    52 = "Parents transfer"
    53 = "Parents transfer"
    54 = "Parents transfer"
    55 = "Parents transfer"
    56 = "Work salary"
    57 = "Parents transfer"
    58 = "Parents transfer"
    59 = "Parents transfer"
    60 = "Security deposit refund"
    61 = "Parents transfer"
    62 = "Parents transfer"
    63 = "Parents transfer"
    64 = "Work salary"
    65 = "Parents transfer"
    66 = "Parents transfer"
    67 = "Parents transfer"
    68 = "Unknown"
    69 = "Parents transfer"
    70 = "Unknown"
    71 = "Work salary"
    72 = "Parents transfer"
    73 = "Parents transfer"
    74 = "Parents transfer"
    75 = "Parents transfer"
    76 = "Work salary"
    77 = "Parents transfer"
    78 = "Unknown"
    79 = "Parents transfer"
    80 = "Parents transfer"
    81 = "Parents transfer"
    82 = "Parents transfer"
    83 = "Scholarship"
    84 = "Parents transfer"
    85 = "Scholarship"
    86 = "Parents transfer"
    87 = "Parents transfer"
    88 = "Grandparents transfer"
    89 = "Scholarship"
    90 = "Grandparents transfer"
    91 = "Tax refund"
    92 = "Parents transfer"
}

foreach ($row in $categoryByRow.Keys) {
    $ws.Cells.Item([int]$row, 3).Value = $categoryByRow[$row]
}

# Narrow the "Income category" column now that the labels are shorter.
$ws.Columns.Item(3).ColumnWidth = 11.498697916666666

# Leave the selection/scroll position where the author left it: bottom of
# the table, category column.
$ws.Range("C91").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 79
$aw.ScrollColumn = 1
